$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Final")
$all = $ws.Cells.FormatConditions
Write-Host $all.Count
for ($i=1; $i -le $all.Count; $i++) {
    $fc = $all.Item($i)
    Write-Host $i $fc.AppliesTo.Address()
}
